# Commit: "added missing outgroup species"
#
# - Adds a new shared string "TODO add 3 outgroup species" to the
#   "Tidied" worksheet, one row below the existing last row of data.
# - Leaves a trail of selection/scroll state matching what a user would
#   end up with after clicking around "All" and then making this edit on
#   "Tidied" (clears the old scrolled-down view on "All", and leaves the
#   selection sitting just below the freshly-typed row on "Tidied").

$wb = $excel.ActiveWorkbook

# Touch the "All" sheet first (clears its old scrolled-down topLeftCell
# state and parks the selection at F37), then come back to "Tidied" so it
# remains the active tab, matching the saved workbook's activeTab.
$wsAll = $wb.Worksheets.Item("All")
$wsAll.Activate()
$wsAll.Range("F37").Select()

# "Tidied" is the sheet that gets the actual content edit: one new row
# appended right after the last existing row (98), naming the missing
# outgroup species entry.
$wsTidied = $wb.Worksheets.Item("Tidied")
$wsTidied.Activate()
$wsTidied.Range("A102").Value = "TODO add 3 outgroup species"
$wsTidied.Range("A103").Select()
